$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
try {
    $nm.Delete()
    Write-Output "deleted OK"
} catch {
    Write-Output "ERROR: $_"
}
